# Updated list of Technologies
# Moves four existing "technology" textboxes on Slide 2 and adds four new
# ones (Kubernetes, Prometheus, Telegraf, ELK), matching the target OOXML.

function Get-ShapeById($Slide, $TargetId) {
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $sh = $Slide.Shapes.Item($i)
        if ($sh.Id -eq $TargetId) {
            return $sh
        }
    }
    return $null
}

# PowerPoint's Shape.Left/Top/Width/Height are single-precision (float32)
# under the hood, so a naive EMU/12700 assignment can land 1 EMU off the
# intended target (the desired integer EMU sometimes falls between two
# adjacent float32 values). Search nearby candidates - first at whole-EMU
# spacing, then at sub-EMU spacing to hop across float32 rounding steps -
# and keep whichever reproduces the exact EMU value.

function Set-PropEmu($Shape, $PropName, $TargetEmu) {
    $bestV = $TargetEmu / 12700.0
    $bestErr = 999999999

    if ($PropName -eq "Left") { $Shape.Left = $bestV; $back0 = [math]::Round($Shape.Left * 12700) }
    elseif ($PropName -eq "Top") { $Shape.Top = $bestV; $back0 = [math]::Round($Shape.Top * 12700) }
    elseif ($PropName -eq "Width") { $Shape.Width = $bestV; $back0 = [math]::Round($Shape.Width * 12700) }
    else { $Shape.Height = $bestV; $back0 = [math]::Round($Shape.Height * 12700) }
    if ($back0 -eq $TargetEmu) { return }

    for ($k = -8; $k -le 8; $k++) {
        $cand = ($TargetEmu + $k) / 12700.0
        if ($PropName -eq "Left") { $Shape.Left = $cand; $back = [math]::Round($Shape.Left * 12700) }
        elseif ($PropName -eq "Top") { $Shape.Top = $cand; $back = [math]::Round($Shape.Top * 12700) }
        elseif ($PropName -eq "Width") { $Shape.Width = $cand; $back = [math]::Round($Shape.Width * 12700) }
        else { $Shape.Height = $cand; $back = [math]::Round($Shape.Height * 12700) }
        $err = [math]::Abs($back - $TargetEmu)
        if ($err -lt $bestErr) {
            $bestErr = $err
            $bestV = $cand
        }
        if ($err -eq 0) {
            if ($PropName -eq "Left") { $Shape.Left = $bestV }
            elseif ($PropName -eq "Top") { $Shape.Top = $bestV }
            elseif ($PropName -eq "Width") { $Shape.Width = $bestV }
            else { $Shape.Height = $bestV }
            return
        }
    }

    for ($step = 0.25; $step -le 4; $step = $step * 2) {
        for ($k = -80; $k -le 80; $k++) {
            $cand = $TargetEmu / 12700.0 + ($k * $step / 12700.0)
            if ($PropName -eq "Left") { $Shape.Left = $cand; $back = [math]::Round($Shape.Left * 12700) }
            elseif ($PropName -eq "Top") { $Shape.Top = $cand; $back = [math]::Round($Shape.Top * 12700) }
            elseif ($PropName -eq "Width") { $Shape.Width = $cand; $back = [math]::Round($Shape.Width * 12700) }
            else { $Shape.Height = $cand; $back = [math]::Round($Shape.Height * 12700) }
            $err = [math]::Abs($back - $TargetEmu)
            if ($err -lt $bestErr) {
                $bestErr = $err
                $bestV = $cand
            }
            if ($err -eq 0) {
                if ($PropName -eq "Left") { $Shape.Left = $bestV }
                elseif ($PropName -eq "Top") { $Shape.Top = $bestV }
                elseif ($PropName -eq "Width") { $Shape.Width = $bestV }
                else { $Shape.Height = $bestV }
                return
            }
        }
    }

    if ($PropName -eq "Left") { $Shape.Left = $bestV }
    elseif ($PropName -eq "Top") { $Shape.Top = $bestV }
    elseif ($PropName -eq "Width") { $Shape.Width = $bestV }
    else { $Shape.Height = $bestV }
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Reposition the four existing textboxes (size unchanged) ---

$xen = Get-ShapeById $s 6
Set-PropEmu $xen "Left" 2895570
Set-PropEmu $xen "Top" 2616536

$packer = Get-ShapeById $s 9
Set-PropEmu $packer "Left" 1671024
Set-PropEmu $packer "Top" 3381299

$dsc = Get-ShapeById $s 12
Set-PropEmu $dsc "Left" 7781369
Set-PropEmu $dsc "Top" 3420368

$poshspec = Get-ShapeById $s 13
Set-PropEmu $poshspec "Left" 4522504
Set-PropEmu $poshspec "Top" 3637515

# --- Add four new textboxes ---
# The slide's shape-id gap (ids 4 & 5 are free from earlier deletions) gets
# filled before new ids continue past the current max (13). Burn through the
# gap with throwaway duplicates first so the real additions land on ids
# 14-17, matching the target file.

$fillerA = $xen.Duplicate()
$fillerB = $xen.Duplicate()
$fillerA.Delete()
$fillerB.Delete()

# Kubernetes (TextBox 13 / id 14) - plain run, no spell-check flag. Height
# is left untouched: a freshly duplicated spAutoFit textbox already reports
# the correct laid-out height (461665 EMU) once the text is changed, and
# nudging it manually only risks landing 1 EMU off via float32 rounding.
$kubernetes = $xen.Duplicate()
$kubernetes.Name = "TextBox 13"
$kubernetes.TextFrame.TextRange.Text = "Kubernetes"
Set-PropEmu $kubernetes "Left" 2323397
Set-PropEmu $kubernetes "Top" 4156558
Set-PropEmu $kubernetes "Width" 1675908

# Prometheus (TextBox 14 / id 15) - plain run, no spell-check flag.
$prometheus = $xen.Duplicate()
$prometheus.Name = "TextBox 14"
$prometheus.TextFrame.TextRange.Text = "Prometheus"
Set-PropEmu $prometheus "Left" 7062747
Set-PropEmu $prometheus "Top" 2606039
Set-PropEmu $prometheus "Width" 1751249

# Telegraf (TextBox 15 / id 16) - flagged err="1" plus endParaRPr, like
# PoshSpec, so duplicate that shape to keep the same run/endParaRPr shape.
$telegraf = $poshspec.Duplicate()
$telegraf.Name = "TextBox 15"
$telegraf.TextFrame.TextRange.Text = "Telegraf"
Set-PropEmu $telegraf "Left" 952974
Set-PropEmu $telegraf "Top" 2729372
Set-PropEmu $telegraf "Width" 1253677

# ELK (TextBox 16 / id 17) - plain run, no spell-check flag.
$elk = $xen.Duplicate()
$elk.Name = "TextBox 16"
$elk.TextFrame.TextRange.Text = "ELK"
Set-PropEmu $elk "Left" 6502192
Set-PropEmu $elk "Top" 3847320
Set-PropEmu $elk "Width" 668773

Write-Host "Done. Slide 2 shape count: $($s.Shapes.Count)"
